# semana 29 de 2025
# Adds week 29 ("AF") column of data to the weekly IRA extract sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header cell (row 1) -------------------------------------------------
# Mirrors the other week-number header cells (e.g. AE1 = "28"), which are
# stored as text. Use a leading apostrophe so the numeric-looking value is
# kept as text instead of being coerced to a number.
$ws.Range("AF1").Value = "'29"

# --- Weekly case counts for column AF (week 29) --------------------------
$afValues = @{
    2  = 65
    3  = 61
    4  = 0
    5  = 1
    6  = 70
    7  = 42
    8  = 32
    9  = 4
    12 = 5
    14 = 3
    15 = 2
    17 = 3
    23 = 7
    24 = 1
    25 = 53
    26 = 3
    27 = 1
    28 = 8
    29 = 0
    30 = 28
    31 = 4
    32 = 7
    34 = 0
    35 = 35
    36 = 4
    37 = 7
    38 = 100
    39 = 1
    40 = 3
    41 = 6
    42 = 26
    43 = 190
    44 = 89
    45 = 156
    46 = 12
    47 = 121
    48 = 6
    49 = 0
    50 = 11
    52 = 46
    53 = 0
    54 = 0
    55 = 9
    56 = 74
    57 = 19
}

foreach ($row in $afValues.Keys) {
    $ws.Cells.Item($row, 32).Value = $afValues[$row]
}

# --- Corrections to existing weeks ---------------------------------------
$ws.Range("AD28").Value = 126
$ws.Range("AE38").Value = 93

# --- New facility name added for row 54 -----------------------------------
$ws.Range("C54").Value = "CLINICA MEDICA TURIN SAS"
